# Careers in Business Intelligence.pptx - apply author's edits:
#  1. Clear the speaker notes text on slide 1 (notes placeholder wiped out).
#  2. Re-cache the auto date field ("datetimeFigureOut") from 3/2/2020 to
#     3/7/2020 everywhere it is cached: the Notes Master, the Slide Master,
#     and every slide layout that carries a Date Placeholder.

$p = $ppt.ActivePresentation

# --- 1. Clear notes text on slide 1 -----------------------------------
$s1 = $p.Slides.Item(1)
$notesPage = $s1.NotesPage
for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
    $shp = $notesPage.Shapes.Item($i)
    if ($shp.Name -like "Notes Placeholder*") {
        $shp.TextFrame.TextRange.Text = ""
    }
}

# --- 2. Refresh the cached "datetimeFigureOut" field text -------------
$newDate = "3/7/2020"

# Notes Master: the Date Placeholder only accepts edits through the
# HeadersFooters.DateAndTime object on this deck.
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = $newDate

# Slide Master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout that has its own Date Placeholder
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
